$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $result = $rng.Find.Execute(
        $find, $false, $false, $false, $false, $false, $true, 1, $false,
        $replace, 2)
    if (-not $result) {
        Write-Output "WARNING: replace failed for: $find"
    }
}

# 1) Fire wall "looks" bullet -> electric fence description
Replace-Text "looks: like a wall of fire (color changes depending on system/level)" `
    "looks: like an electric fence.(number of wires increase as the level/difficulty increases)"

# 2) Fire wall "used for" bullet: buy -> by
Replace-Text "used for: blocking low level attacks powered buy anti-virus software nodes" `
    "used for: blocking low level attacks powered by anti-virus software nodes"

# 3) Anti-virus software nodes "used for" bullet: add repair clause
Replace-Text "used for: power nodes that increase the HP of defenses in area" `
    ("used for: power nodes that increase the HP of defenses in area. " + `
     "can also repair damaged defenses if it isn" + [char]0x2019 + "t dealt with promptly (ex: repair broken sections of the firewall).")

# 4) Remove the whole "Gaser" entry (4 paragraphs of content + 1 spacer paragraph)
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Gaser") {
        $startPara = $p
        $endPara = $d.Paragraphs.Item($i + 4)
        $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
        $rng.Delete()
        break
    }
}

# 5) Spiker "Attack" bullet: under nether -> underneath
Replace-Text "Attack: Shoots a large spike into the ground which travels underground till it shoots up under nether target virus. " `
    "Attack: Shoots a large spike into the ground which travels underground till it shoots up underneath target virus. "

# 6) Level types: over come -> overcome
Replace-Text ": Different types of levels each with unique challenges to over come. " `
    ": Different types of levels each with unique challenges to overcome. "

# 7) Challenges intro: stale/plan. -> stale/plane.
Replace-Text "Used to break up game play to keep it from becoming stale/plan." `
    "Used to break up game play to keep it from becoming stale/plane."

# 8) Budget challenge: spelling fixes
Replace-Text " Budget chalenge: Compleate a set number of levels with limited resources" `
    " Budget challenge: Complete a set number of levels with limited resources"

# 9) Unit challenge: spelling fixes
Replace-Text "Unit chalenge: Compleate a set number of levels with a limited number of units" `
    "Unit challenge: Complete a set number of levels with a limited number of units"

# 10) Boss challenge: spelling fixes
Replace-Text "Boss chalenge: Defeate a rediculously strong enemy that is not in the main game" `
    "Boss challenge: Defeat a ridiculously strong enemy that is not in the main game"

Write-Output "done"
